# Generate Report for Handback
# Updates the status of the db74cb5f-5c6c-4252-91e4-b7c9155840cf.md file
# (row 6 in each sheet) from "Ready for handoff" / stale-handback error
# state to a successfully handed-back state.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E6").Value = "Handed back: in sync with en-US"
$ws.Range("F6").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C6").Value = "Handed back: in sync with en-US"
$ws.Range("L6").Value = "2017-02-17 08:41:29"
$ws.Range("R6").Value = ""

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C6").Value = "Handed back: in sync with en-US"
$ws.Range("L6").Value = "2017-02-17 08:41:52"
$ws.Range("R6").Value = ""
